$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new weekly price rows right after the existing row 250 ---
# (pushes the old rows 251.. down by 2)
$ws.Rows("251:252").Insert()

$ws.Cells.Item(251, 1).Value = 5
$ws.Cells.Item(251, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(251, 3).Value = "Maule"
$ws.Cells.Item(251, 4).Value = 44748
$ws.Cells.Item(251, 5).Value = 7
$ws.Cells.Item(251, 6).Value = "Fruta"
$ws.Cells.Item(251, 7).Value = 100101
$ws.Cells.Item(251, 8).Value = "Berries"
$ws.Cells.Item(251, 9).Value = 100101007
$ws.Cells.Item(251, 10).Value = "Kiwi"
$ws.Cells.Item(251, 11).Value = "Hayward"
$ws.Cells.Item(251, 12).Value = "Especial"
$ws.Cells.Item(251, 13).Value = 190
$ws.Cells.Item(251, 14).Value = 8000
$ws.Cells.Item(251, 15).Value = 8000
$ws.Cells.Item(251, 16).Value = 8000
$ws.Cells.Item(251, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(251, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(251, 19).Value = 444
$ws.Cells.Item(251, 20).Value = 18

$ws.Cells.Item(252, 1).Value = 5
$ws.Cells.Item(252, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(252, 3).Value = "Maule"
$ws.Cells.Item(252, 4).Value = 44748
$ws.Cells.Item(252, 5).Value = 7
$ws.Cells.Item(252, 6).Value = "Fruta"
$ws.Cells.Item(252, 7).Value = 100101
$ws.Cells.Item(252, 8).Value = "Berries"
$ws.Cells.Item(252, 9).Value = 100101007
$ws.Cells.Item(252, 10).Value = "Kiwi"
$ws.Cells.Item(252, 11).Value = "Hayward"
$ws.Cells.Item(252, 12).Value = "Primera"
$ws.Cells.Item(252, 13).Value = 200
$ws.Cells.Item(252, 14).Value = 6000
$ws.Cells.Item(252, 15).Value = 6000
$ws.Cells.Item(252, 16).Value = 6000
$ws.Cells.Item(252, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(252, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(252, 19).Value = 333
$ws.Cells.Item(252, 20).Value = 18

# --- Insert one more new weekly price row further down the table ---
# (original row 268 is now at row 270; push it and everything after down by 1)
$ws.Rows("270:270").Insert()

$ws.Cells.Item(270, 1).Value = 5
$ws.Cells.Item(270, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(270, 3).Value = "Maule"
$ws.Cells.Item(270, 4).Value = 44747
$ws.Cells.Item(270, 5).Value = 7
$ws.Cells.Item(270, 6).Value = "Fruta"
$ws.Cells.Item(270, 7).Value = 100101
$ws.Cells.Item(270, 8).Value = "Berries"
$ws.Cells.Item(270, 9).Value = 100101007
$ws.Cells.Item(270, 10).Value = "Kiwi"
$ws.Cells.Item(270, 11).Value = "Hayward"
$ws.Cells.Item(270, 12).Value = "Primera"
$ws.Cells.Item(270, 13).Value = 300
$ws.Cells.Item(270, 14).Value = 6000
$ws.Cells.Item(270, 15).Value = 6000
$ws.Cells.Item(270, 16).Value = 6000
$ws.Cells.Item(270, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(270, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(270, 19).Value = 333
$ws.Cells.Item(270, 20).Value = 18
